$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C239").Value = 672
$ws.Range("C241").Value = 757
$ws.Range("C243").Value = 375
$ws.Range("C266").Value = 281
$ws.Range("E329").Value = 13
$ws.Range("F329").Value = 12
$ws.Range("G329").Value = 126
$ws.Range("E330").Value = 12
$ws.Range("F330").Value = 9
$ws.Range("G330").Value = 112
$ws.Range("E331").Value = 13
$ws.Range("F331").Value = 9
$ws.Range("G331").Value = 113
$ws.Range("E332").Value = 14
$ws.Range("F332").Value = 8
$ws.Range("G332").Value = 114
$ws.Range("E333").Value = 10
$ws.Range("F333").Value = 8
$ws.Range("G333").Value = 112
$ws.Range("C334").Value = 71
$ws.Range("E334").Value = 11
$ws.Range("F334").Value = 8
$ws.Range("G334").Value = 125
$ws.Range("C335").Value = 152
$ws.Range("E335").Value = 13
$ws.Range("F335").Value = 11
$ws.Range("G335").Value = 133
$ws.Range("E336").Value = 13
$ws.Range("F336").Value = 9
$ws.Range("C337").Value = 92
$ws.Range("E337").Value = 15
$ws.Range("F337").Value = 10
$ws.Range("M337").Value = 2
$ws.Range("C338").Value = 124
$ws.Range("E338").Value = 15
$ws.Range("F338").Value = 12
$ws.Range("G338").Value = 109
$ws.Range("M338").Value = 3
$ws.Range("C339").Value = 107
$ws.Range("G339").Value = 104
$ws.Range("L339").Value = 3
$ws.Range("C340").Value = 76
$ws.Range("E340").Value = 12
$ws.Range("F340").Value = 11
$ws.Range("G340").Value = 103
$ws.Range("L340").Value = 1
$ws.Range("M340").Value = 0
$ws.Range("C341").Value = 36
$ws.Range("E341").Value = 13
$ws.Range("F341").Value = 10
$ws.Range("G341").Value = 112
$ws.Range("L341").Value = 1
$ws.Range("M341").Value = 0
$ws.Range("C342").Value = 16
$ws.Range("E342").Value = 15
$ws.Range("F342").Value = 10
$ws.Range("G342").Value = 116
$ws.Range("L342").Value = 1
$ws.Range("M342").Value = 0
